# Swap the values between row 5 and row 6 for the columns that changed:
# A, B, D, E, F, G, H, Z, AB
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Z", "AB")

foreach ($col in $cols) {
    $addr5 = "$col" + "5"
    $addr6 = "$col" + "6"

    $val5 = $ws.Range($addr5).Value()
    $val6 = $ws.Range($addr6).Value()

    $ws.Range($addr5).Value = $val6
    $ws.Range($addr6).Value = $val5
}
